$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B94 to be a numeric value instead of text
$ws.Range("B94").Value = 5

# Add new row 95 data
$ws.Range("A95").Value = "Ying Tang"
$ws.Range("B95").Value = "'3"
$ws.Range("B95").ClearFormats()
$ws.Range("C95").Value = "it seems like "
$ws.Range("D95").Value = "DFT"
$ws.Range("E95").Value = "WRI"
$ws.Range("F95").Value = "5f16c598-6baf-4325-b317-afc92fb937cc"
$ws.Range("G95").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H95").Value = "Further, it seems like there are errors due to the switching."
